# Add a new "CensusTract" column (G) to the Data sheet: a cleaned/formatted
# version of the raw 6-digit tract code already held in column F
# (e.g. "000107" -> "1.07").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell G1 --------------------------------------------------
# Match the look of the other header cells (bold, centered, bordered) by
# copying F1's formatting onto G1, then set its text.
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "CensusTract"

# --- Data cells G2:G30 ------------------------------------------------
# These are formatted tract numbers such as "1.07", some with a trailing
# space (cleaned-but-not-quite data from the source). Force the range to
# text first so values like "1.07" / "38.03" aren't auto-converted to
# numbers, then restore the (unstyled) default cell formatting.
$ws.Range("G2:G30").NumberFormat = "@"

$ws.Range("G2").Value = "1.07"
$ws.Range("G3").Value = "1.09 "
$ws.Range("G4").Value = "1.15"
$ws.Range("G5").Value = "1.18"
$ws.Range("G6").Value = "1.20"
$ws.Range("G7").Value = "1.21"
$ws.Range("G8").Value = "1.22"
$ws.Range("G9").Value = "1.23"
$ws.Range("G10").Value = "1.25"
$ws.Range("G11").Value = "1.26"
$ws.Range("G12").Value = "1.27"
$ws.Range("G13").Value = "1.28 "
$ws.Range("G14").Value = "1.29"
$ws.Range("G15").Value = "1.30"
$ws.Range("G16").Value = "1.31"
$ws.Range("G17").Value = "1.32"
$ws.Range("G18").Value = "1.34"
$ws.Range("G19").Value = "1.40"
$ws.Range("G20").Value = "1.41"
$ws.Range("G21").Value = "1.42"
$ws.Range("G22").Value = "1.43"
$ws.Range("G23").Value = "1.44"
$ws.Range("G24").Value = "1.45"
$ws.Range("G25").Value = "1.46"
$ws.Range("G26").Value = "12.04 "
$ws.Range("G27").Value = "38.01 "
$ws.Range("G28").Value = "38.03"
$ws.Range("G29").Value = "38.04"
$ws.Range("G30").Value = "9900 "

$ws.Range("G2:G30").ClearFormats()
